# gen_data_test.xlsx edit
#
# 1. Rename the generic "SheetN" tabs to their real, descriptive names.
# 2. On the generation-technology info sheet, mark the two biomass
#    technologies (Biomass_IGCC / Biomass_IGCC_CCS) as no longer
#    resource-limited (g_is_resource_limited: 1 -> 0) because their
#    dispatch is now bounded by the fuel supply curve instead.
# 3. On the generation-cost sheet, drop the now-unused Coal_ST cost rows
#    (2020 & 2030), shifting the remaining rows up.
# 4. Leave the selection / active-tab state the way the workbook was last
#    saved (gen info active, with the relevant ranges highlighted).

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheets -------------------------------------------------
$wb.Worksheets.Item(1).Name = "gen info"
$wb.Worksheets.Item(2).Name = "gen costs"
$wb.Worksheets.Item(3).Name = "ccs info"
$wb.Worksheets.Item(4).Name = "storage info"
$wb.Worksheets.Item(5).Name = "gen_energy"

$wsGenInfo   = $wb.Worksheets.Item("gen info")
$wsGenCosts  = $wb.Worksheets.Item("gen costs")

# --- 2. Biomass technologies are no longer resource limited --------------
# Row 11 = Biomass_IGCC, row 12 = Biomass_IGCC_CCS; column G = g_is_resource_limited
$wsGenInfo.Range("G11").Value = 0
$wsGenInfo.Range("G12").Value = 0

# --- 3. Remove the Coal_ST rows from the gen costs sheet ------------------
$wsGenCosts.Rows("7:8").Delete()

# --- 4. Restore the sheet selections / active tab -------------------------
$wsGenCosts.Activate() | Out-Null
$wsGenCosts.Range("C10:E10").Select() | Out-Null

$wsGenInfo.Activate() | Out-Null
$wsGenInfo.Range("A11:XFD12").Select() | Out-Null
